{"js": "// Replicates the practice-sheet update: a new date heading and 100 new\n// addition/subtraction problems (5 columns x 20 rows) in the table below it.\nconst dateAfter = \"2023-11-05 Sunday\";\nconst grid = [\n  [\"89-76=\", \"88-4=\", \"34+23=\", \"28+65=\", \"28+54=\"],\n  [\"61+18=\", \"48-7=\", \"79+13=\", \"70-39=\", \"16+3=\"],\n  [\"73-34=\", \"34+14=\", \"22+6=\", \"1+19=\", \"66-20=\"],\n  [\"40+31=\", \"22-15=\", \"49-34=\", \"47+29=\", \"95-63=\"],\n  [\"20+65=\", \"67-24=\", \"71-61=\", \"75-46=\", \"23+12=\"],\n  [\"96-64=\", \"75-33=\", \"17+3=\", \"20-5=\", \"76+15=\"],\n  [\"74+17=\", \"58-54=\", \"49+38=\", \"97-74=\", \"33-20=\"],\n  [\"80-14=\", \"63-39=\", \"96-60=\", \"70-2=\", \"39+27=\"],\n  [\"28-12=\", \"58+17=\", \"32+23=\", \"68-50=\", \"92-46=\"],\n  [\"55+30=\", \"42+31=\", \"42+55=\", \"80+0=\", \"38-4=\"],\n  [\"43+25=\", \"2+49=\", \"33-24=\", \"68-33=\", \"38+39=\"],\n  [\"96-18=\", \"80-79=\", \"60-40=\", \"6+81=\", \"33+42=\"],\n  [\"38+3=\", \"52-46=\", \"63-39=\", \"59+29=\", \"70+6=\"],\n  [\"87+5=\", \"39-24=\", \"36+28=\", \"13+15=\", \"22+47=\"],\n  [\"66+32=\", \"93-77=\", \"4+87=\", \"11+55=\", \"73-32=\"],\n  [\"29-10=\", \"99-67=\", \"39+46=\", \"14+31=\", \"81-36=\"],\n  [\"36+48=\", \"77-11=\", \"29+11=\", \"92-41=\", \"91+5=\"],\n  [\"95-84=\", \"71-42=\", \"19-0=\", \"37+62=\", \"11+70=\"],\n  [\"36+0=\", \"12+21=\", \"72-31=\", \"67-63=\", \"17+1=\"],\n  [\"24+66=\", \"96-48=\", \"48-4=\", \"14+20=\", \"84-55=\"]\n];\n\nconst body = context.document.body;\n\n// Update the date paragraph (the single paragraph before the table).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(dateAfter, Word.InsertLocation.replace);\n\n// Update every cell in the table. Replacing the text of the cell's first\n// paragraph (rather than the cell/range itself) keeps the existing run\n// formatting (font, size, etc.) intact.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellParagraphs = [];\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    const cellBody = table.getCell(r, c).body;\n    cellBody.paragraphs.load(\"items\");\n    cellParagraphs.push(cellBody.paragraphs);\n  }\n}\nawait context.sync();\n\nlet k = 0;\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    cellParagraphs[k].items[0].insertText(grid[r][c], Word.InsertLocation.replace);\n    k++;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2023-11-05 Sunday\"\n\n# Update each arithmetic-problem cell in the practice table, preserving\n# cell formatting by assigning to the cell Range.Text rather than replacing\n# the whole cell.\n$t = $d.Tables.Item(1)\n$values = @(\n    @(\"89-76=\", \"88-4=\", \"34+23=\", \"28+65=\", \"28+54=\"),\n    @(\"61+18=\", \"48-7=\", \"79+13=\", \"70-39=\", \"16+3=\"),\n    @(\"73-34=\", \"34+14=\", \"22+6=\", \"1+19=\", \"66-20=\"),\n    @(\"40+31=\", \"22-15=\", \"49-34=\", \"47+29=\", \"95-63=\"),\n    @(\"20+65=\", \"67-24=\", \"71-61=\", \"75-46=\", \"23+12=\"),\n    @(\"96-64=\", \"75-33=\", \"17+3=\", \"20-5=\", \"76+15=\"),\n    @(\"74+17=\", \"58-54=\", \"49+38=\", \"97-74=\", \"33-20=\"),\n    @(\"80-14=\", \"63-39=\", \"96-60=\", \"70-2=\", \"39+27=\"),\n    @(\"28-12=\", \"58+17=\", \"32+23=\", \"68-50=\", \"92-46=\"),\n    @(\"55+30=\", \"42+31=\", \"42+55=\", \"80+0=\", \"38-4=\"),\n    @(\"43+25=\", \"2+49=\", \"33-24=\", \"68-33=\", \"38+39=\"),\n    @(\"96-18=\", \"80-79=\", \"60-40=\", \"6+81=\", \"33+42=\"),\n    @(\"38+3=\", \"52-46=\", \"63-39=\", \"59+29=\", \"70+6=\"),\n    @(\"87+5=\", \"39-24=\", \"36+28=\", \"13+15=\", \"22+47=\"),\n    @(\"66+32=\", \"93-77=\", \"4+87=\", \"11+55=\", \"73-32=\"),\n    @(\"29-10=\", \"99-67=\", \"39+46=\", \"14+31=\", \"81-36=\"),\n    @(\"36+48=\", \"77-11=\", \"29+11=\", \"92-41=\", \"91+5=\"),\n    @(\"95-84=\", \"71-42=\", \"19-0=\", \"37+62=\", \"11+70=\"),\n    @(\"36+0=\", \"12+21=\", \"72-31=\", \"67-63=\", \"17+1=\"),\n    @(\"24+66=\", \"96-48=\", \"48-4=\", \"14+20=\", \"84-55=\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$r - 1][$c - 1]\n    }\n}\n\n"}
